# Scheduled market-data refresh for the Moogle Profits workbook.
# Pulls the latest Universalis current-average-price snapshot for each
# crafting-job leve sheet and recomputes the derived NQ/HQ price + profit
# columns (H:N) that feed off of it.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 963494.0600000001
$ws.Range("I15").Value = 963494.0600000001
$ws.Range("K15").Value = 2890482.18
$ws.Range("M15").Value = -2890313.18
$ws.Range("H32").Value = 2107.4546
$ws.Range("I32").Value = 1673
$ws.Range("J32").Value = 2355.7144
$ws.Range("K32").Value = 1673
$ws.Range("L32").Value = 2355.7144
$ws.Range("M32").Value = -1347
$ws.Range("N32").Value = -3007.7144
$ws.Range("H33").Value = 2575
$ws.Range("I33").Value = 286.375
$ws.Range("K33").Value = 286.375
$ws.Range("M33").Value = -57.375
$ws.Range("H43").Value = 6722.636
$ws.Range("J43").Value = 3632.6667
$ws.Range("L43").Value = 3632.6667
$ws.Range("N43").Value = -3770.6667
$ws.Range("H70").Value = 2777.0908
$ws.Range("I70").Value = 2999.6667
$ws.Range("J70").Value = 2510
$ws.Range("K70").Value = 8999.000100000001
$ws.Range("L70").Value = 7530
$ws.Range("M70").Value = -8729.000100000001
$ws.Range("N70").Value = -8070
$ws.Range("H73").Value = 2777.0908
$ws.Range("I73").Value = 2999.6667
$ws.Range("J73").Value = 2510
$ws.Range("K73").Value = 8999.000100000001
$ws.Range("L73").Value = 7530
$ws.Range("M73").Value = -8063.000100000001
$ws.Range("N73").Value = -9402
$ws.Range("H80").Value = 694.2
$ws.Range("I80").Value = 511.75
$ws.Range("K80").Value = 1535.25
$ws.Range("M80").Value = -537.25
$ws.Range("H83").Value = 694.2
$ws.Range("I83").Value = 511.75
$ws.Range("K83").Value = 4605.75
$ws.Range("M83").Value = 386.25
$ws.Range("H87").Value = 100443.69
$ws.Range("I87").Value = 32000
$ws.Range("J87").Value = 110221.36
$ws.Range("K87").Value = 32000
$ws.Range("L87").Value = 110221.36
$ws.Range("M87").Value = -30752
$ws.Range("N87").Value = -112717.36
$ws.Range("H90").Value = 100443.69
$ws.Range("I90").Value = 32000
$ws.Range("J90").Value = 110221.36
$ws.Range("K90").Value = 96000
$ws.Range("L90").Value = 330664.08
$ws.Range("M90").Value = -89760
$ws.Range("N90").Value = -343144.08
$ws.Range("H111").Value = 68606.64999999999
$ws.Range("I111").Value = 2705
$ws.Range("J111").Value = 96065.664
$ws.Range("K111").Value = 8115
$ws.Range("L111").Value = 288196.992
$ws.Range("M111").Value = -5048
$ws.Range("N111").Value = -294330.992
$ws.Range("H141").Value = 2884.2856
$ws.Range("I141").Value = 2098.2258
$ws.Range("K141").Value = 6294.6774
$ws.Range("M141").Value = -1114.6774

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3187.4614
$ws.Range("I74").Value = 1604.7646
$ws.Range("J74").Value = 4410.4546
$ws.Range("K74").Value = 1604.7646
$ws.Range("L74").Value = 4410.4546
$ws.Range("M74").Value = -730.7646
$ws.Range("N74").Value = -6158.4546
$ws.Range("H77").Value = 3187.4614
$ws.Range("I77").Value = 1604.7646
$ws.Range("J77").Value = 4410.4546
$ws.Range("K77").Value = 8023.823
$ws.Range("L77").Value = 22052.273
$ws.Range("M77").Value = -3655.823
$ws.Range("N77").Value = -30788.273
$ws.Range("H88").Value = 1686.5
$ws.Range("I88").Value = 1219
$ws.Range("J88").Value = 2154
$ws.Range("K88").Value = 1219
$ws.Range("L88").Value = 2154
$ws.Range("M88").Value = -813
$ws.Range("N88").Value = -2966
$ws.Range("H91").Value = 1686.5
$ws.Range("I91").Value = 1219
$ws.Range("J91").Value = 2154
$ws.Range("K91").Value = 1219
$ws.Range("L91").Value = 2154
$ws.Range("M91").Value = 185
$ws.Range("N91").Value = -4962

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1543.5
$ws.Range("I86").Value = 1293
$ws.Range("J86").Value = 4800
$ws.Range("K86").Value = 1293
$ws.Range("L86").Value = 4800
$ws.Range("M86").Value = -170
$ws.Range("N86").Value = -7046
$ws.Range("H89").Value = 1543.5
$ws.Range("I89").Value = 1293
$ws.Range("J89").Value = 4800
$ws.Range("K89").Value = 6465
$ws.Range("L89").Value = 24000
$ws.Range("M89").Value = -849
$ws.Range("N89").Value = -35232
$ws.Range("H105").Value = 2896.077
$ws.Range("I105").Value = 2868.889
$ws.Range("J105").Value = 2957.25
$ws.Range("K105").Value = 2868.889
$ws.Range("L105").Value = 2957.25
$ws.Range("M105").Value = -1121.889
$ws.Range("N105").Value = -6451.25

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6329.0605
$ws.Range("I31").Value = 2643.7727
$ws.Range("K31").Value = 2643.7727
$ws.Range("M31").Value = -2348.7727
$ws.Range("H34").Value = 6329.0605
$ws.Range("I34").Value = 2643.7727
$ws.Range("K34").Value = 2643.7727
$ws.Range("M34").Value = -2441.7727
$ws.Range("H64").Value = 132823.67
$ws.Range("J64").Value = 132823.67
$ws.Range("L64").Value = 132823.67
$ws.Range("N64").Value = -133319.67
$ws.Range("H67").Value = 132823.67
$ws.Range("J67").Value = 132823.67
$ws.Range("L67").Value = 132823.67
$ws.Range("N67").Value = -134539.67
$ws.Range("H105").Value = 2384.625
$ws.Range("I105").Value = 2525.2856
$ws.Range("K105").Value = 2525.2856
$ws.Range("M105").Value = -778.2856000000002
$ws.Range("H112").Value = 199500
$ws.Range("J112").Value = 199500
$ws.Range("L112").Value = 199500
$ws.Range("N112").Value = -202454
$ws.Range("H138").Value = 94727.86
$ws.Range("J138").Value = 94727.86
$ws.Range("L138").Value = 94727.86
$ws.Range("N138").Value = -105007.86

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 225145.4
$ws.Range("I11").Value = 375165
$ws.Range("K11").Value = 1125495
$ws.Range("M11").Value = -1125355
$ws.Range("H14").Value = 6602
$ws.Range("I14").Value = 6602
$ws.Range("K14").Value = 19806
$ws.Range("M14").Value = -19633
$ws.Range("H113").Value = 649.3333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 649.3333
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1947.9999
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -6287.9999
$ws.Range("H129").Value = 13907348
$ws.Range("I129").Value = 14375.5
$ws.Range("K129").Value = 43126.5
$ws.Range("M129").Value = -38126.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10515.929
$ws.Range("I80").Value = 9123.058999999999
$ws.Range("J80").Value = 12668.546
$ws.Range("K80").Value = 9123.058999999999
$ws.Range("L80").Value = 12668.546
$ws.Range("M80").Value = -8125.058999999999
$ws.Range("N80").Value = -14664.546
$ws.Range("H83").Value = 10515.929
$ws.Range("I83").Value = 9123.058999999999
$ws.Range("J83").Value = 12668.546
$ws.Range("K83").Value = 45615.295
$ws.Range("L83").Value = 63342.73
$ws.Range("M83").Value = -40623.295
$ws.Range("N83").Value = -73326.73000000001
$ws.Range("H122").Value = 4575.522
$ws.Range("I122").Value = 2961.85
$ws.Range("J122").Value = 15333.333
$ws.Range("K122").Value = 8885.549999999999
$ws.Range("L122").Value = 45999.999
$ws.Range("M122").Value = -6435.549999999999
$ws.Range("N122").Value = -50899.999
$ws.Range("H132").Value = 4011.9644
$ws.Range("I132").Value = 2518.7368
$ws.Range("K132").Value = 7556.2104
$ws.Range("M132").Value = -5026.2104

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5946.095
$ws.Range("I40").Value = 5020.1763
$ws.Range("J40").Value = 9881.25
$ws.Range("K40").Value = 5020.1763
$ws.Range("L40").Value = 9881.25
$ws.Range("M40").Value = -4884.1763
$ws.Range("N40").Value = -10153.25
$ws.Range("H46").Value = 3591.158
$ws.Range("I46").Value = 927.875
$ws.Range("K46").Value = 927.875
$ws.Range("M46").Value = -739.875
$ws.Range("H68").Value = 5018
$ws.Range("I68").Value = 5110.5454
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 5110.5454
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -4361.5454
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 5018
$ws.Range("I71").Value = 5110.5454
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 25552.727
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -21808.727
$ws.Range("N71").Value = -27488
$ws.Range("H80").Value = 97382.35000000001
$ws.Range("J80").Value = 97382.35000000001
$ws.Range("L80").Value = 97382.35000000001
$ws.Range("N80").Value = -99628.35000000001
$ws.Range("H82").Value = 1101.2
$ws.Range("I82").Value = 1029.25
$ws.Range("J82").Value = 1389
$ws.Range("K82").Value = 1029.25
$ws.Range("L82").Value = 1389
$ws.Range("M82").Value = -668.25
$ws.Range("N82").Value = -2111
$ws.Range("H83").Value = 97382.35000000001
$ws.Range("J83").Value = 97382.35000000001
$ws.Range("L83").Value = 292147.05
$ws.Range("N83").Value = -303379.05
$ws.Range("H85").Value = 1101.2
$ws.Range("I85").Value = 1029.25
$ws.Range("J85").Value = 1389
$ws.Range("K85").Value = 1029.25
$ws.Range("L85").Value = 1389
$ws.Range("M85").Value = 218.75
$ws.Range("N85").Value = -3885
$ws.Range("H122").Value = 5332.0605
$ws.Range("I122").Value = 4427.3706
$ws.Range("K122").Value = 13282.1118
$ws.Range("M122").Value = -10832.1118
$ws.Range("H132").Value = 5573.7617
$ws.Range("I132").Value = 4169.3887
$ws.Range("J132").Value = 14000
$ws.Range("K132").Value = 12508.1661
$ws.Range("L132").Value = 42000
$ws.Range("M132").Value = -9978.166100000002
$ws.Range("N132").Value = -47060

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6867.7144
$ws.Range("I62").Value = 6879.1665
$ws.Range("J62").Value = 6799
$ws.Range("K62").Value = 6879.1665
$ws.Range("L62").Value = 6799
$ws.Range("M62").Value = -6255.1665
$ws.Range("N62").Value = -8047
$ws.Range("H65").Value = 6867.7144
$ws.Range("I65").Value = 6879.1665
$ws.Range("J65").Value = 6799
$ws.Range("K65").Value = 34395.8325
$ws.Range("L65").Value = 33995
$ws.Range("M65").Value = -31275.8325
$ws.Range("N65").Value = -40235
$ws.Range("H132").Value = 3109.7727
$ws.Range("I132").Value = 2327.1052
$ws.Range("K132").Value = 6981.3156
$ws.Range("M132").Value = -4451.3156

